# Weekly fruit/vegetable price update: insert a new record as row 6,
# pushing the existing rows 6-49 down to 7-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 6 (shifts rows 6:49 -> 7:50,
# and copies formatting from the row above as Excel normally does).
$ws.Rows.Item(6).Insert()

# Make sure the date cell keeps the same number format used throughout
# column D (YYYY-MM-DD HH:MM:SS), matching the style used by the other
# rows (copied automatically from row 5 by Insert, but set explicitly
# to be safe).
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat

# Populate the new row with the new weekly record.
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44761
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112010
$ws.Cells.Item(6, 7).Value = "Achicoria"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 55
$ws.Cells.Item(6, 11).Value = 11000
$ws.Cells.Item(6, 12).Value = 11000
$ws.Cells.Item(6, 13).Value = 11000
$ws.Cells.Item(6, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 611
$ws.Cells.Item(6, 17).Value = 18
$ws.Cells.Item(6, 18).Value = "Hortaliza"
